$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status check timestamp in F1
$ws.Range("F1").Value = "Last status check on: 19.02.2022 01:30"

# Update D5 from text "+0.2" to a numeric value 0.2
$ws.Range("D5").Value = 0.2

# Update E5 from text date string to actual Excel date serial value, matching
# the date style used by the other rows in column E
$ws.Range("E5").Value = 44611.05224537037
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
